# Auto-generated Excel COM-interop script
# Adds "Abschnitt zwei" vocabulary rows (28-49) to Sheet1, matching the
# author commit "vocab updated w/ abs 2".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell values, written in the exact order the words were originally
# entered so the shared-string table comes out in the same sequence. ---
$ws.Range("F28").Value = 'der Unterricht'
$ws.Range("H28").Value = 'lesson'
$ws.Range("F29").Value = 'die Klasse'
$ws.Range("H29").Value = 'class'
$ws.Range("F30").Value = 'das Zimmer'
$ws.Range("H30").Value = 'room'
$ws.Range("G28").Value = '''+e '
$ws.Range("A28").Value = 'beginnen'
$ws.Range("B28").Value = 'to begin'
$ws.Range("N29").Value = 'Guten Morgen!'
$ws.Range("O29").Value = 'good morning'
$ws.Range("N30").Value = 'jetzt '
$ws.Range("O30").Value = 'now'
$ws.Range("F31").Value = 'der Satz'
$ws.Range("G31").Value = 'ä, +e'
$ws.Range("H31").Value = 'sentence'
$ws.Range("A29").Value = 'wiederholen'
$ws.Range("B29").Value = 'to repeat'
$ws.Range("F32").Value = 'das Schulzimmer'
$ws.Range("H32").Value = 'classroom'
$ws.Range("F33").Value = 'der Fußboden'
$ws.Range("G33").Value = 'ö'
$ws.Range("H33").Value = 'floor'
$ws.Range("F34").Value = 'die Decke'
$ws.Range("H34").Value = 'ceiling'
$ws.Range("F35").Value = 'die Wand'
$ws.Range("H35").Value = 'wall'
$ws.Range("A30").Value = 'hängen'
$ws.Range("B30").Value = 'to hang'
$ws.Range("B31").Value = 'to be hanging'
$ws.Range("N31").Value = 'vorn'
$ws.Range("O31").Value = 'front'
$ws.Range("N32").Value = 'hier vorn'
$ws.Range("O32").Value = 'up front'
$ws.Range("N33").Value = 'hier oben'
$ws.Range("O33").Value = 'up here'
$ws.Range("N34").Value = 'hinten'
$ws.Range("O34").Value = 'behind'
$ws.Range("N35").Value = 'dort hinten'
$ws.Range("O35").Value = 'over there'
$ws.Range("N36").Value = 'Rechts'
$ws.Range("O36").Value = 'on the right'
$ws.Range("N37").Value = 'Links'
$ws.Range("O37").Value = 'on the left'
$ws.Range("F36").Value = 'die Tafel'
$ws.Range("H36").Value = 'board'
$ws.Range("F37").Value = 'die Kreide'
$ws.Range("H37").Value = 'chalk/crayon'
$ws.Range("F38").Value = 'der Schwamm'
$ws.Range("H38").Value = 'sponge'
$ws.Range("F39").Value = 'die Landkarte'
$ws.Range("H39").Value = 'map'
$ws.Range("F40").Value = 'die Tür'
$ws.Range("H40").Value = 'door'
$ws.Range("F41").Value = 'das Fenster'
$ws.Range("H41").Value = 'window'
$ws.Range("F42").Value = 'die Lampe'
$ws.Range("H42").Value = 'lamp'
$ws.Range("A32").Value = 'verbessern'
$ws.Range("B32").Value = 'to correct'
$ws.Range("F43").Value = 'der Fehler'
$ws.Range("H43").Value = 'mistake'
$ws.Range("N38").Value = 'sehr'
$ws.Range("O38").Value = 'very'
$ws.Range("A33").Value = 'zeigen'
$ws.Range("B33").Value = 'to show'
$ws.Range("A34").Value = 'heißen'
$ws.Range("B34").Value = 'to call'
$ws.Range("A35").Value = 'machen'
$ws.Range("B35").Value = 'to make/to do'
$ws.Range("A36").Value = 'bilden'
$ws.Range("B36").Value = 'to build'
$ws.Range("K28").Value = 'lang'
$ws.Range("L28").Value = 'long'
$ws.Range("K29").Value = 'kurz'
$ws.Range("L29").Value = 'short'
$ws.Range("K30").Value = 'klein'
$ws.Range("L30").Value = 'small'
$ws.Range("L31").Value = 'large'
$ws.Range("K31").Value = 'groß'
$ws.Range("F44").Value = 'das Gegenteil'
$ws.Range("H44").Value = 'opposite/contrary'
$ws.Range("A37").Value = 'diktieren'
$ws.Range("B37").Value = 'to dictate'
$ws.Range("A38").Value = 'schreiben'
$ws.Range("B38").Value = 'to write'
$ws.Range("A39").Value = 'erklären'
$ws.Range("B39").Value = 'to explain'
$ws.Range("A40").Value = 'verstehen'
$ws.Range("B40").Value = 'to understand'
$ws.Range("N39").Value = 'sondern'
$ws.Range("O39").Value = 'instead'
$ws.Range("F45").Value = 'das Kind'
$ws.Range("G45").Value = '''+er'
$ws.Range("H45").Value = 'kid'
$ws.Range("F46").Value = 'das Wort'
$ws.Range("G46").Value = 'ö, +er'
$ws.Range("H46").Value = 'word'
$ws.Range("F47").Value = 'die Stunde'
$ws.Range("H47").Value = 'hour'
$ws.Range("F48").Value = 'die Minute'
$ws.Range("H48").Value = 'minute'
$ws.Range("F49").Value = 'die Sekund'
$ws.Range("H49").Value = 'second'
$ws.Range("A42").Value = 'dauern'
$ws.Range("B42").Value = 'to last'
$ws.Range("N40").Value = 'ist aus'
$ws.Range("O40").Value = 'is over'
$ws.Range("A43").Value = 'schleißen'
$ws.Range("B43").Value = 'to close'
$ws.Range("C28").Value = 'kaufen'
$ws.Range("N28").Value = 'auch'
$ws.Range("O28").Value = 'also'
$ws.Range("C29").Value = 'kaufen'
$ws.Range("G29").Value = '''+n'
$ws.Range("C30").Value = 'kaufen'
$ws.Range("G30").Value = '-'
$ws.Range("C33").Value = 'kaufen'
$ws.Range("C34").Value = 'kaufen'
$ws.Range("G34").Value = '''+n'
$ws.Range("C35").Value = 'kaufen'
$ws.Range("C36").Value = 'antworten'
$ws.Range("G36").Value = '''+n'
$ws.Range("C37").Value = 'kaufen'
$ws.Range("G37").Value = '''+n'
$ws.Range("C38").Value = 'kaufen'
$ws.Range("C39").Value = 'kaufen'
$ws.Range("G39").Value = '''+n'
$ws.Range("C40").Value = 'kaufen'
$ws.Range("G40").Value = '''+en'
$ws.Range("A41").Value = 'lernen'
$ws.Range("B41").Value = 'to learn'
$ws.Range("C41").Value = 'kaufen'
$ws.Range("G41").Value = '-'
$ws.Range("G42").Value = '''+n'
$ws.Range("C43").Value = 'kaufen'
$ws.Range("G43").Value = '-'
$ws.Range("G44").Value = '''+e'
$ws.Range("G47").Value = '''+n'
$ws.Range("G48").Value = '''+n'
$ws.Range("G49").Value = '''+en'
$ws.Range("G35").Value = 'ä, +e'
$ws.Range("G38").Value = 'ä, +e'
$ws.Range("C32").Value = 'verbessern'
$ws.Range("C42").Value = 'verbessern'

# --- Column F widened to fit the longer German nouns added below ---
$ws.Columns("F").ColumnWidth = 18.3

# --- Selection left where the author was last working ---
$ws.Range("B43").Select()

